$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header "president" in column F, copying the header style from E1
$ws.Range("F1").Value = "president"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill column F with "Carter" for every data row (rows 2 through 92)
$lastRow = 92
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = "Carter"
}
